# Update crypto price table with latest values from the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '62.967.57', '  -0.36%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.464.56', '  -0.51%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.01%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '571.96', '  -0.97%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '147.30', '  +0.21%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.999', '  -0.08%  '),
    @(8, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.530', '  -1.77%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.111', '  -0.99%  '),
    @(10, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.163', '  -0.09%  '),
    @(11, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '5.20', '  -1.61%  '),
    @(12, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.348', '  -1.63%  '),
    @(13, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '28.74', '  +0.11%  '),
    @(14, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000175', '  -2.73%  '),
    @(15, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.900.79', '  -0.92%  '),
    @(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '62.822.03', '  -0.57%  '),
    @(17, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.462.81', '  -0.72%  '),
    @(18, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '7.72', '  -6.42%  '),
    @(19, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '10.77', '  -2.77%  '),
    @(20, 'SuiNetwork', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', '2.26', '  -0.33%  '),
    @(21, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '322.62', '  -2.16%  '),
    @(22, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.15', '  +0.27%  '),
    @(23, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.08%  '),
    @(24, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '10.00', '  +3.03%  '),
    @(25, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '64.92', '  -2.20%  '),
    @(26, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '649.69', '  -3.61%  '),
    @(27, 'WrappedeETH', 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth', '2.580.11', '  -2.06%  '),
    @(28, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0966', '  -3.27%  '),
    @(29, 'Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '0.999', '  +0.20%  '),
    @(30, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.42', '  -2.80%  '),
    @(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '7.88', '  -2.28%  '),
    @(32, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.81', '  -2.90%  '),
    @(33, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.133', '  -0.38%  '),
    @(34, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.999', '  -0.02%  '),
    @(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.50', '  -3.65%  '),
    @(36, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '4.66', '  -2.76%  '),
    @(37, 'RenderToken', 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render', '5.36', '  -2.33%  '),
    @(38, 'PolygonEcosystemToken', 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol', '0.364', '  -2.13%  '),
    @(39, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.55', '  -1.46%  '),
    @(40, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '149.56', '  -1.02%  '),
    @(41, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '2.69', '  -2.20%  '),
    @(42, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.72', '  -2.78%  '),
    @(43, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.0₆0313', '  +0.40%  '),
    @(44, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '0.999', '  -0.03%  '),
    @(45, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '152.98', '  -0.75%  '),
    @(46, 'WhiteBITCoin', 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', '15.42', '  +1.92%  '),
    @(47, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.55', '  -1.78%  '),
    @(48, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.605', '  -0.33%  '),
    @(49, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '20.22', '  -2.24%  '),
    @(50, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0508', '  -1.08%  '),
    @(51, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.0904', '  -1.72%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $bCell.NumberFormat = "@"
    $cCell.NumberFormat = "@"
    $dCell.NumberFormat = "@"
    $eCell.NumberFormat = "@"

    $bCell.Value = $row[1]
    $cCell.Value = $row[2]
    $dCell.Value = $row[3]
    $eCell.Value = $row[4]

    $bCell.Style = "Normal"
    $cCell.Style = "Normal"
    $dCell.Style = "Normal"
    $eCell.Style = "Normal"
}
